$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (ANALOG DEVICES) ai_adoption_score ---
$ws.Range("B2").Value = 0.03101405352916946

# --- Insert a brand-new row at position 3 for TELECOMMUNICATION SYS INC,
#     pushing every row currently at 3..50 down by one (to 4..51) ---
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = 23197
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 138613
$ws.Range("D3").Value = "TELECOMMUNICATION SYS INC"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "'541512"

# --- Scattered ai_adoption_score (column B) corrections at their
#     post-insert row positions ---
$ws.Range("B15").Value = 0.1297716338097014
$ws.Range("B18").Value = 0.3336805201127163
$ws.Range("B20").Value = 0.1588691915647324
$ws.Range("B29").Value = 0.1413746127733945
$ws.Range("B33").Value = 0.628501987066414
$ws.Range("B39").Value = 0.1105707743217683
$ws.Range("B48").Value = 0.3457255015383973

# --- Append a brand-new row 52 for GOOGLE INC ---
$ws.Range("A52").Value = 1652044
$ws.Range("B52").Value = 0.4065383566387364
$ws.Range("C52").Value = 160329
$ws.Range("D52").Value = "GOOGLE INC"
$ws.Range("E52").Value = 0.0297704450786113
$ws.Range("F52").Value = "'334111"
